# indigo/spreadsheetform_guides/project.xlsx
# lists: Use get_field_list_from_json from jsondataferret.utils and make
# outcomes table working
#
# This script reproduces, via Excel COM-interop calls, the edits that were
# applied to the workbook:
#   * a new "Outcome" column (M) is added to the Outcomes table, wired up
#     with the SPREADSHEETFORM:DOWN placeholder strings used by the
#     spreadsheetform importer;
#   * the now-"working" data rows (5-17) get placeholder text so the
#     outcomes table isn't empty;
#   * the remembered cell-selections on both sheets move to reflect where
#     the author was last working.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("General Overview")
$wsOutcomes = $wb.Worksheets.Item("Outcomes")

# ---------------------------------------------------------------------
# Outcomes sheet: extend the table with a new column M ("Outcome" down
# the side, under the existing "Source" header group) and populate the
# first body row with placeholder values so the table renders.
# ---------------------------------------------------------------------

# Pull the formatting for the new header cells (M3/M4) from the existing
# D3:D4 cells (same style group, blue band row 3 + bordered row 4) so the
# new column lines up with the existing header band without minting any
# new cell styles.
$wsOutcomes.Range("D3:D4").Copy() | Out-Null
$wsOutcomes.Range("M3:M4").PasteSpecial(-4122) | Out-Null
$wsOutcomes.Range("M4").Value = "Source"

# Pull the formatting for the new data cells (M5:M17) from the existing
# L5:L17 column (same bordered/shaded "data entry" style) so the new
# column's cells share the existing style instead of creating a new one.
$wsOutcomes.Range("L5:L17").Copy() | Out-Null
$wsOutcomes.Range("M5:M17").PasteSpecial(-4122) | Out-Null

# Row 5 becomes the "down" field definition row for the outcomes table:
# title / definition / source placeholders, with "Not in use yet" filling
# the still-unused middle columns.
$wsOutcomes.Range("A5").Value = "SPREADSHEETFORM:DOWN:outcomes:title"
$wsOutcomes.Range("B5").Value = "SPREADSHEETFORM:DOWN:outcomes:definition"
$wsOutcomes.Range("C5:L5").Value = "Not in use yet"
$wsOutcomes.Range("M5").Value = "SPREADSHEETFORM:DOWN:outcomes:source"

# ---------------------------------------------------------------------
# Remembered selections: Outcomes first, General Overview last, so that
# General Overview (sheet 1) ends up the active tab again, matching the
# original file.
# ---------------------------------------------------------------------

$wsOutcomes.Range("J27").Select() | Out-Null
$wsOverview.Range("A4").Select() | Out-Null
